$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Imei**"
$ws.Range("B1").Value = "Satış Tarihi**"
$ws.Range("C1").Value = "Bayi Kodu**"
$ws.Range("D1").Value = "Bayi Adı**"

$ws.Range("D1").Select()
